$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = '21-06-1763'
$ws.Range("E4").Value = 'Non spécifié'
$ws.Range("E5").Value = 'Québec'
$ws.Range("D9").Value = 'Église Saint-Pierre'
$ws.Range("E9").Value = 'Saint-Jean, Île d’Orléans (Qc)'
$ws.Range("D14").Value = 'Paroisse Saint-André'
$ws.Range("D16").Value = 'Inconnu'
$ws.Range("D21").Value = 'Prieuré Saint-Maixent de Verrines'
$ws.Range("E22").Value = 'Hôtel-Dieu de Québec'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1758'
$ws.Range("D29").Value = 'Saint-André'
$ws.Range("D30").Value = 'Paroisse Saint-Jean'
$ws.Range("D31").Value = 'Non spécifié'
$ws.Range("E31").Value = 'Non spécifié'
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = 'France'
$ws.Range("D37").Value = ""
$ws.Range("D39").Value = 'Non spécifié'
$ws.Range("D46").Value = 'Non spécifié'
$ws.Range("D47").Value = 'Saint-Germain'
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = ""
$ws.Range("D53").Value = ""
$ws.Range("D56").Value = 'Non spécifié'
$ws.Range("D57").Value = 'Non spécifié'
$ws.Range("E57").Value = 'Montréal (Qc)'
$ws.Range("D58").Value = ""
$ws.Range("D62").Value = 'Saint-Porchaire'
$ws.Range("E68").Value = 'L''Ancienne-Lorette (Qc)'
$ws.Range("E73").Value = 'Hôtel-Dieu de Québec'
$ws.Range("D78").Value = 'Non spécifié'
$ws.Range("D82").Value = ""
$ws.Range("D83").Value = ""
$ws.Range("E83").Value = ""
$ws.Range("E86").Value = 'Nouvelle-France'
$ws.Range("D93").Value = 'Non spécifié'
$ws.Range("E93").Value = 'Québec'
$ws.Range("D97").Value = 'Paroisse Notre-Dame'
$ws.Range("D99").Value = 'Données introuvables'
$ws.Range("E100").Value = 'France'
$ws.Range("D112").Value = 'Saint-Liguaire'
$ws.Range("E112").Value = 'Non spécifié'
$ws.Range("E113").Value = 'Not Specified'
$ws.Range("D114").Value = 'Paroisse Sainte-Marie-Madeleine'
$ws.Range("D117").Value = 'Non spécifié'
$ws.Range("E117").Value = 'Québec'
$ws.Range("E122").Value = 'Non spécifié'
$ws.Range("D126").Value = ""
$ws.Range("E126").Value = ""
$ws.Range("D130").Value = 'Non spécifié'
$ws.Range("E130").Value = 'Non spécifié'
$ws.Range("E134").Value = 'Charlesbourg (Qc)'
$ws.Range("E136").Value = ""
$ws.Range("D139").Value = 'temple calviniste'
$ws.Range("E139").Value = 'non spécifié'
$ws.Range("D143").Value = 'Données introuvables'
$ws.Range("E145").Value = ""
$ws.Range("D147").Value = 'Saint-Saturnin'
$ws.Range("E151").Value = 'Île-d’Orléans (Qc), Sainte-Famille'
$ws.Range("E154").Value = 'Neuville'
$ws.Range("D156").Value = 'Non spécifié'
$ws.Range("E156").Value = 'Port-Lajoie (Nouveau-Brunswick)'
$ws.Range("D158").Value = 'Non spécifié'
$ws.Range("D159").Value = 'Saint-Pierre'
$ws.Range("D160").Value = 'Saint-André'
$ws.Range("E162").Value = 'Saint-Pierre, Île-d’Orléans (Qc)'
$ws.Range("D164").Value = 'Notre-Dame'
$ws.Range("E164").Value = 'Champlain'
$ws.Range("D165").Value = ""
$ws.Range("E165").Value = ""
$ws.Range("E169").Value = ""
$ws.Range("E173").Value = 'Saint-Joachim (Qc)'
$ws.Range("D176").Value = 'église Notre-Dame'
$ws.Range("E176").Value = ""
$ws.Range("D177").Value = 'Inconnu'
$ws.Range("E177").Value = 'Québec'
$ws.Range("D178").Value = 'Saint-André'
$ws.Range("D180").Value = 'Non spécifié'
$ws.Range("E180").Value = 'Non spécifié'
$ws.Range("D181").Value = 'Inconnu'
$ws.Range("D185").Value = 'Inconnu'
$ws.Range("E185").Value = 'Décès en 1741 (Lieu inconnu)'
$ws.Range("D186").Value = 'Inconnu'
$ws.Range("E186").Value = 'Inconnu'
